# Applies the "cryptos list" refresh described by the commit:
# "Updated cryptos list on Sat Feb 10 11:24:16 UTC 2024 with GitHub Actions"
#
# For every changed row we overwrite the Coin/Link/Price/Volume(1h) cells
# (columns B-E) with their new values. Column D ("Price") sometimes holds
# values that look like plain numbers (e.g. "0.999", "10.00"); Excel would
# normally auto-convert those to numeric cells, so we briefly force a text
# number format, write the literal string, then restore the default style
# so the cell stays un-styled exactly like the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '47.001.60'
$ws.Range("E2").Value = '  +0.65%  '

# Row 3
$ws.Range("D3").Value = '2.478.74'
$ws.Range("E3").Value = '  +0.37%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.41%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.66%  '

# Row 7
$ws.Range("E7").Value = '  -0.25%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.38%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.28%  '

# Row 11
$ws.Range("E11").Value = '  -1.02%  '

# Row 12
$ws.Range("E12").Value = '  +0.23%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.71%  '

# Row 14
$ws.Range("E14").Value = '  -0.28%  '

# Row 15
$ws.Range("D15").Value = '2.865.14'
$ws.Range("E15").Value = '  +0.33%  '

# Row 16
$ws.Range("D16").Value = '2.477.78'
$ws.Range("E16").Value = '  +1.16%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.844'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.13%  '

# Row 18
$ws.Range("D18").Value = '46.927.22'
$ws.Range("E18").Value = '  +0.69%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.38%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.93%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +15.45%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0929'
$ws.Range("E22").Value = '  -0.74%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.52%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '244.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.85%  '

# Row 25
$ws.Range("E25").Value = '  -0.21%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.49%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.59%  '

# Row 29
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.84%  '

# Row 30
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.78'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.65%  '

# Row 31
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.137'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.17%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.45'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.38%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.07%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.32'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.14%  '

# Row 35
$ws.Range("E35").Value = '  +1.32%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.12%  '

# Row 37
$ws.Range("E37").Value = '  +2.12%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.64'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.75%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.71%  '

# Row 40
$ws.Range("E40").Value = '  -0.18%  '

# Row 41
$ws.Range("E41").Value = '  -0.08%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.71%  '

# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '118.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.77%  '

# Row 44
$ws.Range("E44").Value = '  -0.25%  '

# Row 45
$ws.Range("D45").Value = '1.973.72'
$ws.Range("E45").Value = '  -0.20%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.97%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.11%  '

# Row 48
$ws.Range("E48").Value = '  +1.01%  '

# Row 49
$ws.Range("E49").Value = '  -2.81%  '

# Row 50
$ws.Range("E50").Value = '  -4.67%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.46%  '
